# "Creación de datasets finales"
# - Renames the id column header on Hoja1 from "id_matricula" to
#   "id_alum_matric".
# - Adds a new worksheet "Hoja2" (active sheet) that reshapes the data from
#   Hoja1 into id_alum_matric / id_escuela / id_semestre / n_matriculados
#   columns.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Hoja1: header rename -------------------------------------------------
$ws1.Range("A1").Value = "id_alum_matric"

# Approximate the bestFit column width that Excel computed for the renamed
# (longer) header text in column A.
$ws1.Columns.Item(1).ColumnWidth = 13.21875

# Selection moves off the old tabSelected cell now that Hoja2 becomes active.
$ws1.Range("A1:C1").Select() | Out-Null

# --- Hoja2: new sheet, appended after Hoja1 -------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja2"

$ws2.Columns.Item(1).ColumnWidth = 13.21875
$ws2.Columns.Item(4).ColumnWidth = 13.33203125
$ws2.Columns.Item(5).ColumnWidth = 13.33203125

$ws2.Cells.Item(1, 1).Value = "id_alum_matric"
$ws2.Cells.Item(1, 2).Value = "id_escuela"
$ws2.Cells.Item(1, 3).Value = "id_semestre"
$ws2.Cells.Item(1, 4).Value = "n_matriculados"

$rows = @(
  @(1, 1, 314),
  @(2, 21, 130),
  @(1, 2, 978),
  @(2, 22, 534),
  @(1, 3, 911),
  @(2, 23, 498),
  @(1, 5, 999),
  @(2, 25, 594),
  @(1, 6, 927),
  @(2, 26, 572),
  @(1, 8, 880),
  @(2, 28, 590),
  @(1, 9, 1082),
  @(2, 29, 759)
)

$r = 2
foreach ($row in $rows) {
  $ws2.Cells.Item($r, 2).Value = $row[0]
  $ws2.Cells.Item($r, 3).Value = $row[1]
  $ws2.Cells.Item($r, 4).Value = $row[2]
  $r = $r + 1
}

$ws2.Cells.Item(16, 2).Value = 1
$ws2.Cells.Item(17, 2).Value = 2

$ws2.Range("C16").Select() | Out-Null
